# Updated cryptos list values (price + 1h volume change) per the
# target diff. Column D holds numeric-looking price strings (some
# with multiple "." separators, e.g. "24.444.08") that must stay
# plain text, so we prefix them with a leading apostrophe the same
# way a user typing into Excel would force text entry; Excel strips
# the apostrophe and keeps General number format + string storage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.444.08"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "'1.670.23"
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'312.15"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.3959"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D8").Value = "'0.3925"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("D9").Value = "'52.22"
$ws.Range("E9").Value = "  +6.21%  "
$ws.Range("D10").Value = "'1.395"
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'0.08564"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "'24.56"
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").Value = "'7.288"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "'7.966"
$ws.Range("E15").Value = "  +7.75%  "
$ws.Range("E16").Value = "  +5.32%  "
$ws.Range("D17").Value = "'1.667.43"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "'95.02"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'0.07035"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "'20.64"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").Value = "'6.995"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").Value = "'24.445.52"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "'2.499"
$ws.Range("E25").Value = "  +7.43%  "
$ws.Range("D26").Value = "'3.080"
$ws.Range("E26").Value = "  +15.65%  "
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").Value = "'157.18"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'142.68"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "'5.452"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("D31").Value = "'7.961"
$ws.Range("E31").Value = "  -8.42%  "
$ws.Range("D32").Value = "'2.549"
$ws.Range("E32").Value = "  +5.62%  "
$ws.Range("D33").Value = "'1.846.98"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "'1.062"
$ws.Range("E34").Value = "  +12.76%  "
$ws.Range("D35").Value = "'0.03125"
$ws.Range("E35").Value = "  +8.70%  "
$ws.Range("D36").Value = "'0.08274"
$ws.Range("E36").Value = "  +4.10%  "
$ws.Range("D37").Value = "'6.917"
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").Value = "'11.18"
$ws.Range("E38").Value = "  +13.96%  "
$ws.Range("D39").Value = "'0.2765"
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("D40").Value = "'0.09272"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").Value = "'13.73"
$ws.Range("E42").Value = "  +6.25%  "
$ws.Range("D43").Value = "'1.447"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'16.61"
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("D45").Value = "'0.7120"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").Value = "'2.553"
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("D47").Value = "'4.127"
$ws.Range("D48").Value = "'0.9999"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'0.08448"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'136.81"
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("D51").Value = "'1.271"
$ws.Range("E51").Value = "  +2.12%  "
